$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows above row 32. Excel shifts the existing rows 32-59
# down to 34-61 automatically (this reproduces the diff's row-34..61 block,
# which is an exact, unmodified copy of the former rows 32..59).
$ws.Rows("32:33").Insert()

# Row 32 (new) - a "Primera" quality record dated 2023-05-22 (serial 45068)
$ws.Range("A32").Value2 = 7
$ws.Range("B32").Value2 = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C32").Value2 = "Ñuble"
$ws.Range("D32").Value2 = 45068
$ws.Range("D32").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E32").Value2 = 16
$ws.Range("F32").Value2 = 100112044
$ws.Range("G32").Value2 = "Perejil"
$ws.Range("H32").Value2 = "Sin especificar"
$ws.Range("I32").Value2 = "Primera"
$ws.Range("J32").Value2 = 100
$ws.Range("K32").Value2 = 1200
$ws.Range("L32").Value2 = 1200
$ws.Range("M32").Value2 = 1200
$ws.Range("N32").Value2 = "`$/atado 0,5 a 1 kilo"
$ws.Range("O32").Value2 = "Región del Maule"
$ws.Range("P32").Value2 = 1200
$ws.Range("Q32").Value2 = 1
$ws.Range("R32").Value2 = "Hortaliza"

# Row 33 (new) - a "Segunda" quality record dated 2023-05-22 (serial 45068)
$ws.Range("A33").Value2 = 7
$ws.Range("B33").Value2 = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C33").Value2 = "Ñuble"
$ws.Range("D33").Value2 = 45068
$ws.Range("D33").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E33").Value2 = 16
$ws.Range("F33").Value2 = 100112044
$ws.Range("G33").Value2 = "Perejil"
$ws.Range("H33").Value2 = "Sin especificar"
$ws.Range("I33").Value2 = "Segunda"
$ws.Range("J33").Value2 = 100
$ws.Range("K33").Value2 = 1000
$ws.Range("L33").Value2 = 1000
$ws.Range("M33").Value2 = 1000
$ws.Range("N33").Value2 = "`$/atado 0,5 a 1 kilo"
$ws.Range("O33").Value2 = "Región del Maule"
$ws.Range("P33").Value2 = 1000
$ws.Range("Q33").Value2 = 1
$ws.Range("R33").Value2 = "Hortaliza"
